$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2:G5").Value = "2017-01-03 08:28:44"

$zhcn.Range("E2").Value = "mt"
$zhcn.Range("H2").Value = "2017-01-03 08:28:32"

$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2017-01-03 08:28:32"
$zhcn.Range("L3").Value = "2017-01-03 08:29:05"

$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H4").Value = "2017-01-03 08:28:32"
$zhcn.Range("L4").Value = "2017-01-03 08:29:05"

$zhcn.Range("E5").Value = "mt"
$zhcn.Range("H5").Value = "2017-01-03 08:28:32"

$dede.Range("E2").Value = "mt"
$dede.Range("H2").Value = "2017-01-03 08:28:44"

$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2017-01-03 08:28:44"
$dede.Range("L3").Value = "2017-01-03 08:29:16"

$dede.Range("E4").Value = "mt"
$dede.Range("H4").Value = "2017-01-03 08:28:44"
$dede.Range("L4").Value = "2017-01-03 08:29:16"

$dede.Range("E5").Value = "mt"
$dede.Range("H5").Value = "2017-01-03 08:28:44"
